$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Formula = "=3*3"
$ws.Range("B5").Value = $null
